$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Three sentences had been wrapped in extra runs + <w:proofErr/> markers by
#    Word's spell/grammar checker (spellStart/spellEnd, gramStart/gramEnd).
#    The edit cleans these up into a single plain run per sentence. We
#    reproduce that by doing a literal Find & Replace over the whole visible
#    sentence: Word's replace collapses the matched runs (and drops the now
#    orphaned proofErr markers) into one run carrying the replacement text.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "Evaluation based on KPIs, goal achievement, teamwork, and behavior.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Evaluation based on KPIs, goal achievement, teamwork, and behavior.", 2) | Out-Null

$d.Content.Find.Execute(
    "Access to wellness programs, mental health counseling, and fitness initiatives.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Access to wellness programs, mental health counseling, and fitness initiatives.", 2) | Out-Null

$d.Content.Find.Execute(
    " Typically 2–3 days per week, depending on role and manager approval.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Typically 2–3 days per week, depending on role and manager approval.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Add a new dress-code note at the end of the document: a blank spacer
#    paragraph followed by the new sentence, inserted right after the FAQ's
#    last answer ("Report to HR or Facilities...") and before the document's
#    final (already existing) empty trailing paragraph.
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$faqEnd = $lastPara.Previous()

$r = $faqEnd.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$blankPara = $faqEnd.Next()

$br = $blankPara.Range
$br.Collapse(0)
$br.InsertParagraphAfter()
$dressPara = $blankPara.Next()
$dressPara.Range.Text = "Employees need to wear formal dress in office."

Write-Host "HR policy document updated."
